$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'277.81"
$ws.Range("E2").Value = "'0.88%"
$ws.Range("D3").Value = "'27.27"
$ws.Range("E3").Value = "'2.90%"
$ws.Range("D4").Value = "'4.870"
$ws.Range("E4").Value = "'-0.44%"
$ws.Range("D5").Value = "'0.06427"
$ws.Range("E5").Value = "'1.38%"
$ws.Range("D6").Value = "'6.986"
$ws.Range("E6").Value = "'1.34%"
$ws.Range("D7").Value = "'1.192"
$ws.Range("E7").Value = "'-7.56%"
$ws.Range("D8").Value = "'0.8832"
$ws.Range("E8").Value = "'1.88%"
$ws.Range("D9").Value = "'0.1550"
$ws.Range("E9").Value = "'0.90%"
$ws.Range("D10").Value = "'0.05134"
$ws.Range("E10").Value = "'2.09%"
$ws.Range("D11").Value = "'0.07466"
$ws.Range("E11").Value = "'0.67%"
$ws.Range("D12").Value = "'0.02883"
$ws.Range("E12").Value = "'-1.91%"
$ws.Range("D13").Value = "'0.08986"
$ws.Range("E13").Value = "'-0.67%"
$ws.Range("D14").Value = "'0.001577"
$ws.Range("E14").Value = "'0.15%"
$ws.Range("D15").Value = "'0.0006354"
$ws.Range("E15").Value = "'0.57%"
$ws.Range("D16").Value = "'0.006143"
$ws.Range("E16").Value = "'3.93%"
$ws.Range("D17").Value = "'3.483"
$ws.Range("E17").Value = "'1.03%"
$ws.Range("D18").Value = "'3.312"
$ws.Range("E18").Value = "'0.02%"
$ws.Range("E19").Value = "'0.09%"
$ws.Range("D20").Value = "'0.3184"
$ws.Range("E20").Value = "'2.25%"
$ws.Range("D22").Value = "'3.909"
$ws.Range("E22").Value = "'0.39%"
$ws.Range("D23").Value = "'0.04422"
$ws.Range("E23").Value = "'1.20%"
$ws.Range("D24").Value = "'0.1499"
$ws.Range("E24").Value = "'8.63%"
$ws.Range("D26").Value = "'0.001176"
$ws.Range("E26").Value = "'0.69%"
$ws.Range("D27").Value = "'0.003874"
$ws.Range("E27").Value = "'-8.94%"
$ws.Range("E28").Value = "'-1.56%"
$ws.Range("E29").Value = "'15.74%"
$ws.Range("D40").Value = "'0.04142"
$ws.Range("E40").Value = "'0.81%"
$ws.Range("D41").Value = "'0.006797"
$ws.Range("E41").Value = "'-2.67%"
$ws.Range("E42").Value = "'0.37%"
$ws.Range("D43").Value = "'0.001920"
$ws.Range("E43").Value = "'-10.45%"
$ws.Range("D44").Value = "'0.01140"
$ws.Range("E44").Value = "'5.77%"
$ws.Range("D45").Value = "'0.00005315"
$ws.Range("E45").Value = "'0.93%"
$ws.Range("D46").Value = "'1.687"
$ws.Range("E46").Value = "'13.27%"
$ws.Range("D47").Value = "'0.01851"
$ws.Range("E47").Value = "'-7.32%"

$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Style = "Normal"
